# Swap the data of several row-pairs back to their correct match rows.
# For each pair the "id/grouping" columns A (row #), C (Div) and D (Date)
# stay put, while everything describing the actual match (B and E..AD)
# is exchanged between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to swap: B (2) and E..AD (5..30)
$colsToSwap = @(2) + (5..30)

# Row pairs that need to be swapped (1-based worksheet rows)
$rowPairs = @(
    @(69, 70),
    @(120, 121),
    @(161, 162),
    @(180, 181)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    # Read all values for both rows first so writes don't clobber reads
    $valuesA = @{}
    $valuesB = @{}
    foreach ($col in $colsToSwap) {
        $valuesA[$col] = $ws.Cells.Item($rowA, $col).Value()
        $valuesB[$col] = $ws.Cells.Item($rowB, $col).Value()
    }

    foreach ($col in $colsToSwap) {
        $ws.Cells.Item($rowA, $col).Value = $valuesB[$col]
        $ws.Cells.Item($rowB, $col).Value = $valuesA[$col]
    }
}
